$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Update technician name in H3 (was "Dakota Myers", now "Franz Ferdinand")
$ws.Range("H3").Value = "Franz Ferdinand"

# Update the active selection to match the saved view state (H3)
[void]$ws.Range("H3").Select()
